# Patient_List.xlsx — update the note in column H, row 3.
# The rest of the sheet (patient records in A:F, the other H-column
# notes) is left untouched.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H3").Value = "cock"
